$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh computed TPM-derived values for the Cxcl12-Cxcr4 LR-pair sheet
$ws.Range("G2").Value = 160.9723563333334
$ws.Range("H2").Value = 482.917069
$ws.Range("I2").Value = 0.3931645655589854
$ws.Range("J2").Value = 0.3931645655589854
$ws.Range("M2").Value = 13.71977066666667
$ws.Range("N2").Value = 41.159312
$ws.Range("O2").Value = 0.5515038136402627
$ws.Range("P2").Value = 0.5515038136402626
$ws.Range("Q2").Value = 2208.503812566281
$ws.Range("R2").Value = 19876.53431309653
$ws.Range("S2").Value = 0.2168317572939975
$ws.Range("T2").Value = 0.2168317572939975
$ws.Range("G3").Value = 160.9723563333334
$ws.Range("H3").Value = 482.917069
$ws.Range("I3").Value = 0.3931645655589854
$ws.Range("J3").Value = 0.3931645655589854
$ws.Range("O3").Value = 0.172077867958883
$ws.Range("P3").Value = 0.1720778679588829
$ws.Range("Q3").Value = 689.0879410914841
$ws.Range("R3").Value = 6201.791469823357
$ws.Range("S3").Value = 0.06765492019837067
$ws.Range("T3").Value = 0.06765492019837066
$ws.Range("G4").Value = 160.9723563333334
$ws.Range("H4").Value = 482.917069
$ws.Range("I4").Value = 0.3931645655589854
$ws.Range("J4").Value = 0.3931645655589854
$ws.Range("O4").Value = 0.2764183184008545
$ws.Range("P4").Value = 0.2764183184008545
$ws.Range("Q4").Value = 1106.920559663886
$ws.Range("R4").Value = 9962.285036974979
$ws.Range("S4").Value = 0.1086778880666172
$ws.Range("T4").Value = 0.1086778880666172
$ws.Range("I5").Value = 0.2197635343237224
$ws.Range("J5").Value = 0.2197635343237224
$ws.Range("M5").Value = 13.71977066666667
$ws.Range("N5").Value = 41.159312
$ws.Range("O5").Value = 0.5515038136402627
$ws.Range("P5").Value = 0.5515038136402626
$ws.Range("Q5").Value = 1234.466800757929
$ws.Range("R5").Value = 11110.20120682136
$ws.Range("S5").Value = 0.1212004272785956
$ws.Range("T5").Value = 0.1212004272785956
$ws.Range("I6").Value = 0.2197635343237224
$ws.Range("J6").Value = 0.2197635343237224
$ws.Range("O6").Value = 0.172077867958883
$ws.Range("P6").Value = 0.1720778679588829
$ws.Range("S6").Value = 0.03781644044153494
$ws.Range("T6").Value = 0.03781644044153494
$ws.Range("I7").Value = 0.2197635343237224
$ws.Range("J7").Value = 0.2197635343237224
$ws.Range("O7").Value = 0.2764183184008545
$ws.Range("P7").Value = 0.2764183184008545
$ws.Range("S7").Value = 0.06074666660359179
$ws.Range("T7").Value = 0.06074666660359181
$ws.Range("I8").Value = 0.3870719001172923
$ws.Range("J8").Value = 0.3870719001172923
$ws.Range("M8").Value = 13.71977066666667
$ws.Range("N8").Value = 41.159312
$ws.Range("O8").Value = 0.5515038136402627
$ws.Range("P8").Value = 0.5515038136402626
$ws.Range("Q8").Value = 2174.279785185942
$ws.Range("R8").Value = 19568.51806667347
$ws.Range("S8").Value = 0.2134716290676695
$ws.Range("T8").Value = 0.2134716290676695
$ws.Range("I9").Value = 0.3870719001172923
$ws.Range("J9").Value = 0.3870719001172923
$ws.Range("O9").Value = 0.172077867958883
$ws.Range("P9").Value = 0.1720778679588829
$ws.Range("Q9").Value = 678.4095060219187
$ws.Range("R9").Value = 6105.685554197268
$ws.Range("S9").Value = 0.06660650731897735
$ws.Range("T9").Value = 0.06660650731897735
$ws.Range("I10").Value = 0.3870719001172923
$ws.Range("J10").Value = 0.3870719001172923
$ws.Range("O10").Value = 0.2764183184008545
$ws.Range("P10").Value = 0.2764183184008545
$ws.Range("S10").Value = 0.1069937637306454
$ws.Range("T10").Value = 0.1069937637306454

Write-Output "Updated TPM values on sheet $($ws.Name)"
